$wb = $excel.ActiveWorkbook

# --- Step1_Data ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("AE2").Value = 0.0009379242810787525
$ws.Range("AF2").Value = 0.001698601146577091
$ws.Range("AI2").Value = 0.04467381447538246
$ws.Range("D2").Value = 0.2605598194680626
$ws.Range("F2").Value = 0.1333392950682297
$ws.Range("G2").Value = 0.01557568371055711
$ws.Range("H2").Value = 0.05144960314866772
$ws.Range("I2").Value = 0.001891878845235774
$ws.Range("J2").Value = 0.01102334401284815
$ws.Range("K2").Value = 0.00465182763312846
$ws.Range("L2").Value = 0.002763541283223176
$ws.Range("N2").Value = 0.02075608389873396
$ws.Range("O2").Value = 0.1605308550461695
$ws.Range("P2").Value = 0.03133180736799802
$ws.Range("Q2").Value = 0.1343854541377556
$ws.Range("S2").Value = 0.0809223338361748
$ws.Range("T2").Value = 0.01135239219991543
$ws.Range("U2").Value = 0.02990757964465863
$ws.Range("Z2").Value = 0.002248160795603079
$ws.Range("AB3").Value = 0.004583697710067895
$ws.Range("AC3").Value = 0.005746246653726636
$ws.Range("AD3").Value = 0.0151769952252822
$ws.Range("AI3").Value = 0.0301466239310643
$ws.Range("D3").Value = 0.3134382243586167
$ws.Range("F3").Value = 0.1428383548134479
$ws.Range("G3").Value = 0.01153690470924683
$ws.Range("H3").Value = 0.04432734569614721
$ws.Range("J3").Value = 0.01881961104627377
$ws.Range("L3").Value = 0.004130124788082941
$ws.Range("N3").Value = 0.0096280222331679
$ws.Range("O3").Value = 0.1741062794079408
$ws.Range("P3").Value = 0.07177018864365892
$ws.Range("Q3").Value = 0.04872808538064579
$ws.Range("R3").Value = 0.003341768899028189
$ws.Range("S3").Value = 0.03314731615159337
$ws.Range("T3").Value = 0.00506446169986638
$ws.Range("U3").Value = 0.02213569322971907
$ws.Range("X3").Value = 0.02428874147716304
$ws.Range("Z3").Value = 0.01704531394526032
$ws.Range("AA4").Value = 0.00916444107626488
$ws.Range("AE4").Value = 0.01392814647267863
$ws.Range("AF4").Value = 0.0009503843313881941
$ws.Range("AJ4").Value = 0.02971610109137296
$ws.Range("E4").Value = 0.2371450737113429
$ws.Range("G4").Value = 0.1718116493888919
$ws.Range("H4").Value = 0.01155098345208545
$ws.Range("I4").Value = 0.04795883080954184
$ws.Range("K4").Value = 0.01566704258236823
$ws.Range("L4").Value = 0.02330033608720819
$ws.Range("M4").Value = 0.01147911876082595
$ws.Range("O4").Value = 0.04833891815611251
$ws.Range("P4").Value = 0.1446110901864938
$ws.Range("Q4").Value = 0.06393150330817049
$ws.Range("R4").Value = 0.07085246636581628
$ws.Range("S4").Value = 0.02611413669649911
$ws.Range("T4").Value = 0.04004091201069983
$ws.Range("U4").Value = 0.009376289962775033
$ws.Range("V4").Value = 0.02406257554946404
$ws.Range("AA5").Value = 0.02912870875757922
$ws.Range("AE5").Value = 0.0002384894648214687
$ws.Range("AF5").Value = 0.0008359842936227797
$ws.Range("AG5").Value = 0.004596580553926574
$ws.Range("AH5").Value = 0.01343835554110361
$ws.Range("AJ5").Value = 0.02051784604630679
$ws.Range("E5").Value = 0.2713070837673133
$ws.Range("G5").Value = 0.1681218698637197
$ws.Range("H5").Value = 0.02748009696032527
$ws.Range("I5").Value = 0.05152285502235969
$ws.Range("K5").Value = 0.007275039541223416
$ws.Range("L5").Value = 0.02424184661209161
$ws.Range("M5").Value = 0.003825665148053303
$ws.Range("O5").Value = 0.054935088523334
$ws.Range("P5").Value = 0.1143084733691273
$ws.Range("Q5").Value = 0.02492813296486079
$ws.Range("R5").Value = 0.08431682811602485
$ws.Range("S5").Value = 0.02908413862659222
$ws.Range("T5").Value = 0.03111450661451666
$ws.Range("U5").Value = 0.01322504805457976
$ws.Range("V5").Value = 0.02538503571592458
$ws.Range("Z5").Value = 0.0001723264425930052
$ws.Range("AE6").Value = 0.002893553666939558
$ws.Range("AF6").Value = 0.009942521052731577
$ws.Range("AG6").Value = 0.01032072745644992
$ws.Range("AI6").Value = 0.01650906383773083
$ws.Range("D6").Value = 0.1924368212432876
$ws.Range("F6").Value = 0.09506947079763566
$ws.Range("G6").Value = 0.07249514301565829
$ws.Range("H6").Value = 0.08967607557102369
$ws.Range("N6").Value = 0.1639004799288875
$ws.Range("O6").Value = 0.1702928832913681
$ws.Range("Q6").Value = 0.03460803766118675
$ws.Range("R6").Value = 0.005820300075650163
$ws.Range("S6").Value = 0.04615675445203903
$ws.Range("T6").Value = 0.001701107425599377
$ws.Range("U6").Value = 0.07130827401405458
$ws.Range("Z6").Value = 0.01686878650975714

# --- Step2_Sj ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("AA2").Value = 0.9526896600969617
$ws.Range("AB2").Value = 0.9526896600969617
$ws.Range("AC2").Value = 0.9526896600969617
$ws.Range("AD2").Value = 0.9526896600969617
$ws.Range("AE2").Value = 0.9536275843780404
$ws.Range("AF2").Value = 0.9553261855246175
$ws.Range("AG2").Value = 0.9553261855246175
$ws.Range("AH2").Value = 0.9553261855246175
$ws.Range("AI2").Value = 0.9999999999999999
$ws.Range("AJ2").Value = 0.9999999999999999
$ws.Range("D2").Value = 0.2605598194680626
$ws.Range("E2").Value = 0.2605598194680626
$ws.Range("F2").Value = 0.3938991145362922
$ws.Range("G2").Value = 0.4094747982468493
$ws.Range("H2").Value = 0.460924401395517
$ws.Range("I2").Value = 0.4628162802407528
$ws.Range("J2").Value = 0.473839624253601
$ws.Range("K2").Value = 0.4784914518867294
$ws.Range("L2").Value = 0.4812549931699526
$ws.Range("M2").Value = 0.4812549931699526
$ws.Range("N2").Value = 0.5020110770686866
$ws.Range("O2").Value = 0.6625419321148561
$ws.Range("P2").Value = 0.6938737394828541
$ws.Range("Q2").Value = 0.8282591936206097
$ws.Range("R2").Value = 0.8282591936206097
$ws.Range("S2").Value = 0.9091815274567845
$ws.Range("T2").Value = 0.9205339196567
$ws.Range("U2").Value = 0.9504414993013586
$ws.Range("V2").Value = 0.9504414993013586
$ws.Range("W2").Value = 0.9504414993013586
$ws.Range("X2").Value = 0.9504414993013586
$ws.Range("Y2").Value = 0.9504414993013586
$ws.Range("Z2").Value = 0.9526896600969617
$ws.Range("AA3").Value = 0.9443464364798594
$ws.Range("AB3").Value = 0.9489301341899273
$ws.Range("AC3").Value = 0.9546763808436539
$ws.Range("AD3").Value = 0.9698533760689361
$ws.Range("AE3").Value = 0.9698533760689361
$ws.Range("AF3").Value = 0.9698533760689361
$ws.Range("AG3").Value = 0.9698533760689361
$ws.Range("AH3").Value = 0.9698533760689361
$ws.Range("D3").Value = 0.3134382243586167
$ws.Range("E3").Value = 0.3134382243586167
$ws.Range("F3").Value = 0.4562765791720647
$ws.Range("G3").Value = 0.4678134838813115
$ws.Range("H3").Value = 0.5121408295774588
$ws.Range("I3").Value = 0.5121408295774588
$ws.Range("J3").Value = 0.5309604406237325
$ws.Range("K3").Value = 0.5309604406237325
$ws.Range("L3").Value = 0.5350905654118154
$ws.Range("M3").Value = 0.5350905654118154
$ws.Range("N3").Value = 0.5447185876449834
$ws.Range("O3").Value = 0.7188248670529243
$ws.Range("P3").Value = 0.7905950556965832
$ws.Range("Q3").Value = 0.8393231410772289
$ws.Range("R3").Value = 0.8426649099762571
$ws.Range("S3").Value = 0.8758122261278505
$ws.Range("T3").Value = 0.8808766878277169
$ws.Range("U3").Value = 0.903012381057436
$ws.Range("V3").Value = 0.903012381057436
$ws.Range("W3").Value = 0.903012381057436
$ws.Range("X3").Value = 0.9273011225345991
$ws.Range("Y3").Value = 0.9273011225345991
$ws.Range("Z3").Value = 0.9443464364798594
$ws.Range("AA4").Value = 0.9554053681045603
$ws.Range("AB4").Value = 0.9554053681045603
$ws.Range("AC4").Value = 0.9554053681045603
$ws.Range("AD4").Value = 0.9554053681045603
$ws.Range("AE4").Value = 0.9693335145772389
$ws.Range("AF4").Value = 0.9702838989086271
$ws.Range("AG4").Value = 0.9702838989086271
$ws.Range("AH4").Value = 0.9702838989086271
$ws.Range("AI4").Value = 0.9702838989086271
$ws.Range("E4").Value = 0.2371450737113429
$ws.Range("F4").Value = 0.2371450737113429
$ws.Range("G4").Value = 0.4089567231002348
$ws.Range("H4").Value = 0.4205077065523203
$ws.Range("I4").Value = 0.4684665373618621
$ws.Range("J4").Value = 0.4684665373618621
$ws.Range("K4").Value = 0.4841335799442303
$ws.Range("L4").Value = 0.5074339160314385
$ws.Range("M4").Value = 0.5189130347922645
$ws.Range("N4").Value = 0.5189130347922645
$ws.Range("O4").Value = 0.5672519529483769
$ws.Range("P4").Value = 0.7118630431348707
$ws.Range("Q4").Value = 0.7757945464430411
$ws.Range("R4").Value = 0.8466470128088575
$ws.Range("S4").Value = 0.8727611495053565
$ws.Range("T4").Value = 0.9128020615160564
$ws.Range("U4").Value = 0.9221783514788314
$ws.Range("V4").Value = 0.9462409270282954
$ws.Range("W4").Value = 0.9462409270282954
$ws.Range("X4").Value = 0.9462409270282954
$ws.Range("Y4").Value = 0.9462409270282954
$ws.Range("Z4").Value = 0.9462409270282954
$ws.Range("AA5").Value = 0.9603727441002187
$ws.Range("AB5").Value = 0.9603727441002187
$ws.Range("AC5").Value = 0.9603727441002187
$ws.Range("AD5").Value = 0.9603727441002187
$ws.Range("AE5").Value = 0.9606112335650402
$ws.Range("AF5").Value = 0.9614472178586629
$ws.Range("AG5").Value = 0.9660437984125895
$ws.Range("AH5").Value = 0.9794821539536931
$ws.Range("AI5").Value = 0.9794821539536931
$ws.Range("AJ5").Value = 0.9999999999999999
$ws.Range("E5").Value = 0.2713070837673133
$ws.Range("F5").Value = 0.2713070837673133
$ws.Range("G5").Value = 0.4394289536310331
$ws.Range("H5").Value = 0.4669090505913583
$ws.Range("I5").Value = 0.518431905613718
$ws.Range("J5").Value = 0.518431905613718
$ws.Range("K5").Value = 0.5257069451549414
$ws.Range("L5").Value = 0.549948791767033
$ws.Range("M5").Value = 0.5537744569150863
$ws.Range("N5").Value = 0.5537744569150863
$ws.Range("O5").Value = 0.6087095454384203
$ws.Range("P5").Value = 0.7230180188075476
$ws.Range("Q5").Value = 0.7479461517724084
$ws.Range("R5").Value = 0.8322629798884332
$ws.Range("S5").Value = 0.8613471185150254
$ws.Range("T5").Value = 0.8924616251295421
$ws.Range("U5").Value = 0.9056866731841219
$ws.Range("V5").Value = 0.9310717089000464
$ws.Range("W5").Value = 0.9310717089000464
$ws.Range("X5").Value = 0.9310717089000464
$ws.Range("Y5").Value = 0.9310717089000464
$ws.Range("Z5").Value = 0.9312440353426394
$ws.Range("AA6").Value = 0.9603341339861478
$ws.Range("AB6").Value = 0.9603341339861478
$ws.Range("AC6").Value = 0.9603341339861478
$ws.Range("AD6").Value = 0.9603341339861478
$ws.Range("AE6").Value = 0.9632276876530873
$ws.Range("AF6").Value = 0.9731702087058188
$ws.Range("AG6").Value = 0.9834909361622688
$ws.Range("AH6").Value = 0.9834909361622688
$ws.Range("AI6").Value = 0.9999999999999997
$ws.Range("AJ6").Value = 0.9999999999999997
$ws.Range("D6").Value = 0.1924368212432876
$ws.Range("E6").Value = 0.1924368212432876
$ws.Range("F6").Value = 0.2875062920409233
$ws.Range("G6").Value = 0.3600014350565816
$ws.Range("H6").Value = 0.4496775106276052
$ws.Range("I6").Value = 0.4496775106276052
$ws.Range("J6").Value = 0.4496775106276052
$ws.Range("K6").Value = 0.4496775106276052
$ws.Range("L6").Value = 0.4496775106276052
$ws.Range("M6").Value = 0.4496775106276052
$ws.Range("N6").Value = 0.6135779905564928
$ws.Range("O6").Value = 0.7838708738478608
$ws.Range("P6").Value = 0.7838708738478608
$ws.Range("Q6").Value = 0.8184789115090475
$ws.Range("R6").Value = 0.8242992115846977
$ws.Range("S6").Value = 0.8704559660367367
$ws.Range("T6").Value = 0.8721570734623361
$ws.Range("U6").Value = 0.9434653474763907
$ws.Range("V6").Value = 0.9434653474763907
$ws.Range("W6").Value = 0.9434653474763907
$ws.Range("X6").Value = 0.9434653474763907
$ws.Range("Y6").Value = 0.9434653474763907
$ws.Range("Z6").Value = 0.9603341339861478

# --- Step3_DataPts_0.5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("D2").Value = 13
$ws.Range("F2").Value = 0.5020110770686866
$ws.Range("G2").Value = 12
$ws.Range("D3").Value = 7
$ws.Range("F3").Value = 0.5121408295774588
$ws.Range("G3").Value = 6
$ws.Range("D4").Value = 11
$ws.Range("F4").Value = 0.5074339160314385
$ws.Range("G4").Value = 9
$ws.Range("D5").Value = 8
$ws.Range("F5").Value = 0.518431905613718
$ws.Range("G5").Value = 6
$ws.Range("F6").Value = 0.6135779905564928

# --- Step3_DataPts_0.7 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value = 16
$ws.Range("F2").Value = 0.8282591936206097
$ws.Range("G2").Value = 15
$ws.Range("F3").Value = 0.7188248670529243
$ws.Range("F4").Value = 0.7118630431348707
$ws.Range("D5").Value = 15
$ws.Range("F5").Value = 0.7230180188075476
$ws.Range("G5").Value = 13
$ws.Range("F6").Value = 0.7838708738478608

# --- Step3_DataPts_0.8 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F2").Value = 0.8282591936206097
$ws.Range("D3").Value = 16
$ws.Range("F3").Value = 0.8393231410772289
$ws.Range("G3").Value = 15
$ws.Range("D4").Value = 17
$ws.Range("F4").Value = 0.8466470128088575
$ws.Range("G4").Value = 15
$ws.Range("D5").Value = 17
$ws.Range("F5").Value = 0.8322629798884332
$ws.Range("G5").Value = 15
$ws.Range("D6").Value = 16
$ws.Range("F6").Value = 0.8184789115090475
$ws.Range("G6").Value = 15

# --- Step3_DataPts_0.9 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F2").Value = 0.9091815274567845
$ws.Range("D3").Value = 20
$ws.Range("F3").Value = 0.903012381057436
$ws.Range("G3").Value = 19
$ws.Range("D4").Value = 19
$ws.Range("F4").Value = 0.9128020615160564
$ws.Range("G4").Value = 17
$ws.Range("D5").Value = 20
$ws.Range("F5").Value = 0.9056866731841219
$ws.Range("G5").Value = 18
$ws.Range("D6").Value = 20
$ws.Range("F6").Value = 0.9434653474763907
$ws.Range("G6").Value = 19
